$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 14
$ws.Range("B17").Value = 17
$ws.Range("C17").Value = 1789621
$ws.Range("D17").Value = "M Safroni"
$ws.Range("E17").Value = "Jakarta"
$ws.Range("F17").Value = "2000-07-13"
$ws.Range("G17").Value = "Islam"
$ws.Range("H17").Value = "Laki-laki"
$ws.Range("I17").Value = "A"
$ws.Range("J17").Value = "Belum Nikah"
$ws.Range("K17").Value = "PTT"
$ws.Range("L17").Value = "jakarta pusat"
$ws.Range("M17").Value = "082180712764"
$ws.Range("N17").Value = "apalah@gmail.com"
$ws.Range("S17").Value = "Cetak"
$ws.Range("T17").Value = "2058-07-13"
